$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update H2:H11 values from 50 to 60
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 60
}

# Update the active selection shown in the sheet view from H16 to H19
$ws.Range("H19").Select()
